$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 currently exists as an empty formatted placeholder row (same visual
# band as rows 2-9). Bring over the same cell formatting (borders/number
# formats/fills) used by the row directly above it, then fill in the new
# task's data.
$ws.Range("A9:J9").Copy()
$ws.Range("A10:J10").PasteSpecial(-4122)

$ws.Range("A10").Value = "adminpage"
$ws.Range("B10").Value = 45289
$ws.Range("C10").Value = 45326
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "thanh"

# Move the active selection like the author left it.
$ws.Range("F17").Select()
